# Add a new "Label" column (H) to the worksheet: 0 for Control rows, 1 for MDD rows.
# (This mirrors re-running the classification export with an added ground-truth
# diagnosis label column, as used for the manuscript figure.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell H1 ---------------------------------------------------------
$ws.Range("H1").Value = "Label"

# Match the header formatting used by the other header cells (bold, centered,
# bordered) by copying the format from the neighboring header cell G1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows ---------------------------------------------------------------
# Column A in each block cycles through the five Control patients followed by
# the five MDD patients (rows 2-11, then repeated in rows 12-21 for a second
# "Iterations" setting). The new Label column records the binary diagnosis
# label: 0 = Control, 1 = MDD.
$labels = @(0, 0, 0, 0, 0, 1, 1, 1, 1, 1)

$blockStarts = @(2, 12)
foreach ($start in $blockStarts) {
    for ($i = 0; $i -lt $labels.Length; $i++) {
        $row = $start + $i
        $ws.Cells.Item($row, 8).Value = $labels[$i]
    }
}
